$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Martin Guptill'
$ws.Range("D2").Value = 'LBW'
$ws.Range("E2").Value = ' Tabraiz Shamsi'
$ws.Range("J2").Value = 'Temba Bavuma(C)'
$ws.Range("K2").Value = 7
$ws.Range("N2").Value = ' Tim Southee'
$ws.Range("A3").Value = 'Daryl Mitchell'
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 2
$ws.Range("E3").Value = ' Anrich Nortje'
$ws.Range("J3").Value = 'Quinton de Kock'
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("N3").Value = ' Trent Boult'
$ws.Range("A4").Value = 'Kane Williamson(C)'
$ws.Range("B4").Value = 46
$ws.Range("C4").Value = 18
$ws.Range("D4").Value = 'NOT OUT'
$ws.Range("E4").Value = ' '
$ws.Range("J4").Value = 'Rassie Va der Dussen'
$ws.Range("K4").Value = 28
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 'LBW'
$ws.Range("N4").Value = ' Mitchell Santner'
$ws.Range("A5").Value = 'Devon Conway'
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 'Bowled'
$ws.Range("E5").Value = ' Anrich Nortje'
$ws.Range("J5").Value = 'Aiden Markram'
$ws.Range("K5").Value = 9
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 'Bowled'
$ws.Range("N5").Value = ' Ish Sodhi'
$ws.Range("A6").Value = 'Glenn Phillips'
$ws.Range("B6").Value = 16
$ws.Range("C6").Value = 5
$ws.Range("E6").Value = ' Kagiso Rabada'
$ws.Range("J6").Value = 'David Miller'
$ws.Range("K6").Value = 43
$ws.Range("L6").Value = 17
$ws.Range("M6").Value = 'NOT OUT'
$ws.Range("N6").Value = ' '
$ws.Range("A7").Value = 'James Neesham'
$ws.Range("B7").Value = 16
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 'LBW'
$ws.Range("E7").Value = ' Tabraiz Shamsi'
$ws.Range("J7").Value = 'Reeza Hendricks'
$ws.Range("K7").Value = 24
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = '* NOT OUT'
$ws.Range("N7").Value = ' '
$ws.Range("A8").Value = 'Mitchell Santner'
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 'LBW'
$ws.Range("E8").Value = ' Anrich Nortje'
$ws.Range("J8").Value = 'Dwaine Pretorius'
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = ' '
$ws.Range("N8").Value = ' '
$ws.Range("A9").Value = 'Adam Milne'
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("E9").Value = ' Keshav Maharaj'
$ws.Range("J9").Value = 'Kagiso Rabada'
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = ' '
$ws.Range("N9").Value = ' '
$ws.Range("A10").Value = 'Ish Sodhi'
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 'Caught'
$ws.Range("E10").Value = ' Keshav Maharaj'
$ws.Range("J10").Value = 'Keshav Maharaj'
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ' '
$ws.Range("N10").Value = ' '
$ws.Range("A11").Value = 'Tim Southee'
$ws.Range("B11").Value = 0
$ws.Range("E11").Value = ' Keshav Maharaj'
$ws.Range("J11").Value = 'Anrich Nortje'
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ' '
$ws.Range("N11").Value = ' '
$ws.Range("A12").Value = 'Trent Boult'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 'Bowled'
$ws.Range("E12").Value = ' Keshav Maharaj'
$ws.Range("J12").Value = 'Tabraiz Shamsi'
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ' '
$ws.Range("N12").Value = ' '
$ws.Range("A16").Value = 111
$ws.Range("B16").Value = 10
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '8.0'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 48
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = 4
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = '7.1'
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = 43
$ws.Range("A21").Value = 'Kagiso Rabada'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '1.0'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 16
$ws.Range("J21").Value = 'Adam Milne'
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = '1.0'
$ws.Range("K21").Style = "Normal"
$ws.Range("L21").Value = 23
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 23
$ws.Range("A22").Value = 'Dwaine Pretorius'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '1.0'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = 23
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 23
$ws.Range("J22").Value = 'Mitchell Santner'
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = '1.0'
$ws.Range("K22").Style = "Normal"
$ws.Range("L22").Value = 15
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 15
$ws.Range("A23").Value = 'Tabraiz Shamsi'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '2.0'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 9.5
$ws.Range("J23").Value = 'Trent Boult'
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = '2.0'
$ws.Range("K23").Style = "Normal"
$ws.Range("L23").Value = 28
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 14
$ws.Range("A24").Value = 'Anrich Nortje'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '2.0'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 13
$ws.Range("J24").Value = 'Tim Southee'
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = '2.0'
$ws.Range("K24").Style = "Normal"
$ws.Range("L24").Value = 27
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 13.5
$ws.Range("A25").Value = 'Keshav Maharaj'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2.0'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 13.5
$ws.Range("J25").Value = 'Ish Sodhi'
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = '1.1'
$ws.Range("K25").Style = "Normal"
$ws.Range("L25").Value = 20
$ws.Range("M25").Value = 1
$ws.Range("N25").Value = 18.18
